$d = $word.ActiveDocument

# --- Change 1: Title -- insert a new "Dummy " run before the existing "Trial " run ---
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertBefore("Dummy ")

# --- Changes 2-5: walk the paragraphs once, matching by paraId-bearing text, and
#     fix up the four paragraphs that changed -- use the live index so earlier
#     edits (which do not add/remove paragraphs) can't desync us. ---
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text

    if ($i -eq 12 -and $t -like "Input:*abcd@xyz*") {
        # TC002 "Input: " + "abcd@xyz" (spell-check split runs) -> single run, no proofErr
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="11D19FE3" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000"><w:r><w:t>Input: abcd@xyz</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
    }
    elseif ($i -eq 20 -and $t -like "Input:*abcd@xyz*") {
        # TC004 "Input: " + "abcd@xyz" -> single run, no proofErr, gains the
        # lastRenderedPageBreak that used to sit on the following paragraph
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="27E0FDDE" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000"><w:r><w:lastRenderedPageBreak/><w:t>Input: abcd@xyz</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
    }
    elseif ($i -eq 21 -and $t -like "Expected Result: Rejected*Must include at least one number*") {
        # TC004 "Expected Result" paragraph loses its lastRenderedPageBreak
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7D780310" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000"><w:r><w:t>Expected Result: Rejected – Must include at least one number</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
    }
    elseif ($i -eq 36 -and $t -like "Password:*abcdefgh*") {
        # Input 2 "Password: " + "abcdefgh" (spell-check split runs) -> single run
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1A472765" w14:textId="77777777" w:rsidR="002E4A0C" w:rsidRDefault="00000000"><w:r><w:t>Password: abcdefgh</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
    }
}
